$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.693535480230485
$ws.Range("D2").Value = 8.964373503298315
$ws.Range("E2").Value = 13.60319223818741
$ws.Range("F2").Value = 33.86260390055643
$ws.Range("G2").Value = 3.662661072502026
$ws.Range("I2").Value = 21.50744172468364
$ws.Range("J2").Value = 10.1476926252566
$ws.Range("K2").Value = 11.4962342157431
$ws.Range("M2").Value = 15.78047382490055
$ws.Range("O2").Value = 25.61809616040857
$ws.Range("B3").Value = 5.551646108867483
$ws.Range("D3").Value = 8.941986685789756
$ws.Range("E3").Value = 13.6145532777809
$ws.Range("F3").Value = 33.94177352349425
$ws.Range("G3").Value = 3.664681031797967
$ws.Range("I3").Value = 21.62941994279307
$ws.Range("J3").Value = 10.17346321834326
$ws.Range("K3").Value = 11.15744222379996
$ws.Range("M3").Value = 15.65253548411965
$ws.Range("O3").Value = 25.70657796104332
$ws.Range("B4").Value = 5.463329917961605
$ws.Range("D4").Value = 8.929483008116652
$ws.Range("E4").Value = 13.62379522648399
$ws.Range("F4").Value = 33.99867488474119
$ws.Range("G4").Value = 3.665987382031759
$ws.Range("I4").Value = 21.70809870867007
$ws.Range("J4").Value = 10.19047468736195
$ws.Range("K4").Value = 10.9450671579013
$ws.Range("M4").Value = 15.57556599518092
$ws.Range("O4").Value = 25.76666578012498
$ws.Range("B5").Value = 5.427091471654792
$ws.Range("D5").Value = 8.924703425089596
$ws.Range("E5").Value = 13.62813165509174
$ws.Range("F5").Value = 34.02394237503192
$ws.Range("G5").Value = 3.666536400916627
$ws.Range("I5").Value = 21.74111514158511
$ws.Range("J5").Value = 10.19770611806712
$ws.Range("K5").Value = 10.85755086863554
$ws.Range("M5").Value = 15.54462533148404
$ws.Range("O5").Value = 25.7925970247012
$ws.Range("B6").Value = 5.421060636122607
$ws.Range("D6").Value = 8.923928947341144
$ws.Range("E6").Value = 13.62888616672574
$ws.Range("F6").Value = 34.02826347408138
$ws.Range("G6").Value = 3.666628573487996
$ws.Range("I6").Value = 21.74665520772959
$ws.Range("J6").Value = 10.19892496743933
$ws.Range("K6").Value = 10.84296394663816
$ws.Range("M6").Value = 15.5395140865019
$ws.Range("O6").Value = 25.79699006114219
$ws.Range("B7").Value = 5.462842131031601
$ws.Range("D7").Value = 8.929417266041202
$ws.Range("E7").Value = 13.62385139969849
$ws.Range("F7").Value = 33.99900723735692
$ws.Range("G7").Value = 3.66599471872454
$ws.Range("I7").Value = 21.70854011264725
$ws.Range("J7").Value = 10.19057100140353
$ws.Range("K7").Value = 10.94389064747662
$ws.Range("M7").Value = 15.57514696328182
$ws.Range("O7").Value = 25.76700965184025
$ws.Range("B8").Value = 5.644890386688367
$ws.Range("D8").Value = 8.956398924043949
$ws.Range("E8").Value = 13.60663945100803
$ws.Range("F8").Value = 33.88817789992537
$ws.Range("G8").Value = 3.663343869709307
$ws.Range("I8").Value = 21.54871619015558
$ws.Range("J8").Value = 10.15633191002043
$ws.Range("K8").Value = 11.38039403259184
$ws.Range("M8").Value = 15.73604774257071
$ws.Range("O8").Value = 25.64740729933164
$ws.Range("B9").Value = 5.990265796133508
$ws.Range("D9").Value = 9.019001674579611
$ws.Range("E9").Value = 13.59084604709929
$ws.Range("F9").Value = 33.73683387258228
$ws.Range("G9").Value = 3.658667541130416
$ws.Range("I9").Value = 21.26519604560319
$ws.Range("J9").Value = 10.09860233919464
$ws.Range("K9").Value = 12.19684904986466
$ws.Range("M9").Value = 16.0629613896243
$ws.Range("O9").Value = 25.45871281702875
$ws.Range("B10").Value = 6.23428014887828
$ws.Range("D10").Value = 9.070673327476852
$ws.Range("E10").Value = 13.5901538918395
$ws.Range("F10").Value = 33.66612513160338
$ws.Range("G10").Value = 3.655546704282436
$ws.Range("I10").Value = 21.07493800258813
$ws.Range("J10").Value = 10.0619066413668
$ws.Range("K10").Value = 12.76659484336105
$ws.Range("M10").Value = 16.30846070454198
$ws.Range("O10").Value = 25.34822012423446
$ws.Range("B11").Value = 6.342677664455702
$ws.Range("D11").Value = 9.095359049608382
$ws.Range("E11").Value = 13.59219709495648
$ws.Range("F11").Value = 33.64279084373427
$ws.Range("G11").Value = 3.654194611453383
$ws.Range("I11").Value = 20.99226446931828
$ws.Range("J11").Value = 10.0464502131401
$ws.Range("K11").Value = 13.01816204300965
$ws.Range("M11").Value = 16.42094642998227
$ws.Range("O11").Value = 25.30410345382186
$ws.Range("B12").Value = 6.383313128357155
$ws.Range("D12").Value = 9.10487165910447
$ws.Range("E12").Value = 13.59330862057985
$ws.Range("F12").Value = 33.63522709144529
$ws.Range("G12").Value = 3.653692274668997
$ws.Range("I12").Value = 20.96151262747671
$ws.Range("J12").Value = 10.04077477437624
$ws.Range("K12").Value = 13.11225029214109
$ws.Range("M12").Value = 16.46362802382045
$ws.Range("O12").Value = 25.28828452976761
$ws.Range("B13").Value = 6.374580457275311
$ws.Range("D13").Value = 9.102815707169311
$ws.Range("E13").Value = 13.59305423175556
$ws.Range("F13").Value = 33.63679945653659
$ws.Range("G13").Value = 3.653800032506749
$ws.Range("I13").Value = 20.96811095023659
$ws.Range("J13").Value = 10.04198918792513
$ws.Range("K13").Value = 13.0920401400557
$ws.Range("M13").Value = 16.45443244987799
$ws.Range("O13").Value = 25.2916519152942
$ws.Range("B14").Value = 6.346029236177484
$ws.Range("D14").Value = 9.096138387126954
$ws.Range("E14").Value = 13.59228177859806
$ws.Range("F14").Value = 33.642143055942
$ws.Range("G14").Value = 3.654153090347936
$ws.Range("I14").Value = 20.98972339141335
$ws.Range("J14").Value = 10.04597973393393
$ws.Range("K14").Value = 13.02592671351503
$ws.Range("M14").Value = 16.42445632648439
$ws.Range("O14").Value = 25.30278422501145
$ws.Range("B15").Value = 6.328486030548221
$ws.Range("D15").Value = 9.092069627155256
$ws.Range("E15").Value = 13.59185258045707
$ws.Range("F15").Value = 33.64558193580066
$ws.Range("G15").Value = 3.654370606552859
$ws.Range("I15").Value = 21.00303381708144
$ws.Range("J15").Value = 10.04844717671733
$ws.Range("K15").Value = 12.98527510434186
$ws.Range("M15").Value = 16.40610534585368
$ws.Range("O15").Value = 25.30971870423536
$ws.Range("B16").Value = 6.227140392554345
$ws.Range("D16").Value = 9.069083377788377
$ws.Range("E16").Value = 13.59006770884952
$ws.Range("F16").Value = 33.6678280213446
$ws.Range("G16").Value = 3.655636422717038
$ws.Range("I16").Value = 21.08041868887102
$ws.Range("J16").Value = 10.06294161796625
$ws.Range("K16").Value = 12.74999410859434
$ws.Range("M16").Value = 16.30112313953028
$ws.Range("O16").Value = 25.35122722053592
$ws.Range("B17").Value = 6.164272935548452
$ws.Range("D17").Value = 9.055280781751215
$ws.Range("E17").Value = 13.58957585256209
$ws.Range("F17").Value = 33.68373920305311
$ws.Range("G17").Value = 3.656430237358117
$ws.Range("I17").Value = 21.12888269835488
$ws.Range("J17").Value = 10.07215003951903
$ws.Range("K17").Value = 12.60364734314463
$ws.Range("M17").Value = 16.23690501601984
$ws.Range("O17").Value = 25.37826805325893
$ws.Range("B18").Value = 6.127870006803815
$ws.Range("D18").Value = 9.047453277317901
$ws.Range("E18").Value = 13.58951499009707
$ws.Range("F18").Value = 33.69372204727787
$ws.Range("G18").Value = 3.656893183430594
$ws.Range("I18").Value = 21.15712287220011
$ws.Range("J18").Value = 10.07756289222843
$ws.Range("K18").Value = 12.5187602202982
$ws.Range("M18").Value = 16.20004621082096
$ws.Range("O18").Value = 25.39439963102439
$ws.Range("B19").Value = 6.115504023851939
$ws.Range("D19").Value = 9.0448222988311
$ws.Range("E19").Value = 13.58953255477932
$ws.Range("F19").Value = 33.69724473171269
$ws.Range("G19").Value = 3.657051023665107
$ws.Range("I19").Value = 21.16674726832602
$ws.Range("J19").Value = 10.07941559272862
$ws.Range("K19").Value = 12.48989913120707
$ws.Range("M19").Value = 16.18758073327399
$ws.Range("O19").Value = 25.3999607487018
$ws.Range("B20").Value = 6.170990739631479
$ws.Range("D20").Value = 9.056738599624165
$ws.Range("E20").Value = 13.58960523987002
$ws.Range("F20").Value = 33.68195938922694
$ws.Range("G20").Value = 3.656345076042767
$ws.Range("I20").Value = 21.12368586829589
$ws.Range("J20").Value = 10.0711577409209
$ws.Range("K20").Value = 12.61930053630758
$ws.Range("M20").Value = 16.24373331685013
$ws.Range("O20").Value = 25.3753296241108
$ws.Range("B21").Value = 6.354426896146251
$ws.Range("D21").Value = 9.098095250028255
$ws.Range("E21").Value = 13.5924995098473
$ws.Range("F21").Value = 33.64053896127029
$ws.Range("G21").Value = 3.654049126604411
$ws.Range("I21").Value = 20.9833602554478
$ws.Range("J21").Value = 10.04480279648633
$ws.Range("K21").Value = 13.04537830986529
$ws.Range("M21").Value = 16.43325895627638
$ws.Range("O21").Value = 25.29949029566718
$ws.Range("B22").Value = 6.471894480887584
$ws.Range("D22").Value = 9.126081520957396
$ws.Range("E22").Value = 13.59635937290542
$ws.Range("F22").Value = 33.62088572191448
$ws.Range("G22").Value = 3.652604941317956
$ws.Range("I22").Value = 20.89488218823331
$ws.Range("J22").Value = 10.0286132459287
$ws.Range("K22").Value = 13.31696580242797
$ws.Range("M22").Value = 16.55761040309911
$ws.Range("O22").Value = 25.25509697042627
$ws.Range("B23").Value = 6.409432722706453
$ws.Range("D23").Value = 9.111058838375282
$ws.Range("E23").Value = 13.59411966096681
$ws.Range("F23").Value = 33.63069568725893
$ws.Range("G23").Value = 3.653370589929462
$ws.Range("I23").Value = 20.94180961537741
$ws.Range("J23").Value = 10.03715929988948
$ws.Range("K23").Value = 13.17266820036942
$ws.Range("M23").Value = 16.49120726808377
$ws.Range("O23").Value = 25.27831627538746
$ws.Range("B24").Value = 6.167954426093065
$ws.Range("D24").Value = 9.056079184037525
$ws.Range("E24").Value = 13.58959126254697
$ws.Range("F24").Value = 33.68276144134354
$ws.Range("G24").Value = 3.656383557000105
$ws.Range("I24").Value = 21.12603417886307
$ws.Range("J24").Value = 10.07160598901795
$ws.Range("K24").Value = 12.61222605889117
$ws.Range("M24").Value = 16.24064604956223
$ws.Range("O24").Value = 25.37665626452248
$ws.Range("B25").Value = 5.89836224862373
$ws.Range("D25").Value = 9.001051313024126
$ws.Range("E25").Value = 13.5931990992034
$ws.Range("F25").Value = 33.77068308598585
$ws.Range("G25").Value = 3.659877078276183
$ws.Range("I25").Value = 21.33871408039054
$ws.Range("J25").Value = 10.11321408120795
$ws.Range("K25").Value = 12.19684904986466
$ws.Range("M25").Value = 15.97346983561924
$ws.Range("O25").Value = 25.50483121501016
